$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H64").Value = 3900.0386
$ws.Range("I64").Value = 3635.6365
$ws.Range("K64").Value = 3635.6365
$ws.Range("M64").Value = -3387.6365

$ws.Range("H67").Value = 3900.0386
$ws.Range("I67").Value = 3635.6365
$ws.Range("K67").Value = 3635.6365
$ws.Range("M67").Value = -2777.6365

$ws.Range("H86").Value = 3162.375
$ws.Range("I86").Value = 1899.5
$ws.Range("K86").Value = 1899.5
$ws.Range("M86").Value = -776.5

$ws.Range("H89").Value = 3162.375
$ws.Range("I89").Value = 1899.5
$ws.Range("K89").Value = 9497.5
$ws.Range("M89").Value = -3881.5

$ws.Range("H138").Value = 2979.3538
$ws.Range("I138").Value = 1291.3513
$ws.Range("J138").Value = 5209.9287
$ws.Range("K138").Value = 3874.0539
$ws.Range("L138").Value = 15629.7861
$ws.Range("M138").Value = 1265.9461
$ws.Range("N138").Value = -25909.7861

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 21173.172
$ws.Range("I32").Value = 23086.154
$ws.Range("J32").Value = 12883.583
$ws.Range("K32").Value = 23086.154
$ws.Range("L32").Value = 12883.583
$ws.Range("M32").Value = -22799.154
$ws.Range("N32").Value = -13457.583

$ws.Range("H61").Value = 7974.396
$ws.Range("I61").Value = 4507.7666
$ws.Range("J61").Value = 13752.111
$ws.Range("K61").Value = 4507.7666
$ws.Range("L61").Value = 13752.111
$ws.Range("M61").Value = -4295.7666
$ws.Range("N61").Value = -14176.111

$ws.Range("H74").Value = 6042.36
$ws.Range("I74").Value = 2141.682
$ws.Range("J74").Value = 34647.332
$ws.Range("K74").Value = 2141.682
$ws.Range("L74").Value = 34647.332
$ws.Range("M74").Value = -1267.682
$ws.Range("N74").Value = -36395.332

$ws.Range("H77").Value = 6042.36
$ws.Range("I77").Value = 2141.682
$ws.Range("J77").Value = 34647.332
$ws.Range("K77").Value = 10708.41
$ws.Range("L77").Value = 173236.66
$ws.Range("M77").Value = -6340.41
$ws.Range("N77").Value = -181972.66

$ws.Range("H132").Value = 1749.8193
$ws.Range("I132").Value = 1339.2858
$ws.Range("J132").Value = 3043
$ws.Range("K132").Value = 4017.8574
$ws.Range("L132").Value = 9129
$ws.Range("M132").Value = -1487.8574
$ws.Range("N132").Value = -14189

$ws.Range("H136").Value = 7974.396
$ws.Range("I136").Value = 4507.7666
$ws.Range("J136").Value = 13752.111
$ws.Range("K136").Value = 13523.2998
$ws.Range("L136").Value = 41256.333
$ws.Range("M136").Value = -10973.2998
$ws.Range("N136").Value = -46356.333

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H55").Value = 60000
$ws.Range("J55").Value = 60000
$ws.Range("L55").Value = 60000
$ws.Range("N55").Value = -60546

$ws.Range("H86").Value = 1588.9778
$ws.Range("I86").Value = 1605.5897
$ws.Range("J86").Value = 1481
$ws.Range("K86").Value = 1605.5897
$ws.Range("L86").Value = 1481
$ws.Range("M86").Value = -482.5897
$ws.Range("N86").Value = -3727

$ws.Range("H89").Value = 1588.9778
$ws.Range("I89").Value = 1605.5897
$ws.Range("J89").Value = 1481
$ws.Range("K89").Value = 8027.9485
$ws.Range("L89").Value = 7405
$ws.Range("M89").Value = -2411.9485
$ws.Range("N89").Value = -18637

$ws.Range("H99").Value = 1049.875
$ws.Range("I99").Value = 1000
$ws.Range("J99").Value = 1399
$ws.Range("K99").Value = 1000
$ws.Range("L99").Value = 1399
$ws.Range("M99").Value = 498
$ws.Range("N99").Value = -4395

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 2245.5715
$ws.Range("I31").Value = 1588.6111
$ws.Range("J31").Value = 4064.8462
$ws.Range("K31").Value = 1588.6111
$ws.Range("L31").Value = 4064.8462
$ws.Range("M31").Value = -1293.6111
$ws.Range("N31").Value = -4654.8462

$ws.Range("H34").Value = 2245.5715
$ws.Range("I34").Value = 1588.6111
$ws.Range("J34").Value = 4064.8462
$ws.Range("K34").Value = 1588.6111
$ws.Range("L34").Value = 4064.8462
$ws.Range("M34").Value = -1386.6111
$ws.Range("N34").Value = -4468.8462

$ws.Range("H58").Value = 1820171
$ws.Range("I58").Value = 2458286.5
$ws.Range("J58").Value = 3996
$ws.Range("K58").Value = 2458286.5
$ws.Range("L58").Value = 3996
$ws.Range("M58").Value = -2458083.5
$ws.Range("N58").Value = -4402

$ws.Range("H62").Value = 3637.0715
$ws.Range("I62").Value = 3339.5
$ws.Range("K62").Value = 3339.5
$ws.Range("M62").Value = -2715.5

$ws.Range("H65").Value = 3637.0715
$ws.Range("I65").Value = 3339.5
$ws.Range("K65").Value = 16697.5
$ws.Range("M65").Value = -13577.5

$ws.Range("H132").Value = 3305.3484
$ws.Range("I132").Value = 3595.3696
$ws.Range("J132").Value = 2638.3
$ws.Range("K132").Value = 10786.1088
$ws.Range("L132").Value = 7914.900000000001
$ws.Range("M132").Value = -8256.1088
$ws.Range("N132").Value = -12974.9

$ws.Range("H134").Value = 2488.6035
$ws.Range("I134").Value = 1676.5312
$ws.Range("J134").Value = 3488.077
$ws.Range("K134").Value = 5029.5936
$ws.Range("L134").Value = 10464.231
$ws.Range("M134").Value = -2494.5936
$ws.Range("N134").Value = -15534.231

$ws.Range("H136").Value = 1820171
$ws.Range("I136").Value = 2458286.5
$ws.Range("J136").Value = 3996
$ws.Range("K136").Value = 7374859.5
$ws.Range("L136").Value = 11988
$ws.Range("M136").Value = -7372309.5
$ws.Range("N136").Value = -17088

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H113").Value = 2789.647
$ws.Range("I113").Value = 4200
$ws.Range("J113").Value = 2355.6924
$ws.Range("K113").Value = 4200
$ws.Range("L113").Value = 2355.6924
$ws.Range("M113").Value = -2030
$ws.Range("N113").Value = -6695.6924

$ws.Range("H126").Value = 2633.88
$ws.Range("I126").Value = 1773.1818
$ws.Range("J126").Value = 3310.1428
$ws.Range("K126").Value = 5319.5454
$ws.Range("L126").Value = 9930.428400000001
$ws.Range("M126").Value = -2849.5454
$ws.Range("N126").Value = -14870.4284

$ws.Range("H132").Value = 3981.8044
$ws.Range("I132").Value = 1758.8684
$ws.Range("J132").Value = 14540.75
$ws.Range("K132").Value = 5276.6052
$ws.Range("L132").Value = 43622.25
$ws.Range("M132").Value = -2746.6052
$ws.Range("N132").Value = -48682.25

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H132").Value = 5120.1953
$ws.Range("I132").Value = 5699
$ws.Range("J132").Value = 4303.0586
$ws.Range("K132").Value = 17097
$ws.Range("L132").Value = 12909.1758
$ws.Range("M132").Value = -14567
$ws.Range("N132").Value = -17969.1758

$ws.Range("H136").Value = 3243.838
$ws.Range("I136").Value = 1760.9183
$ws.Range("J136").Value = 6150.36
$ws.Range("K136").Value = 5282.7549
$ws.Range("L136").Value = 18451.08
$ws.Range("M136").Value = -2732.7549
$ws.Range("N136").Value = -23551.08

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H62").Value = 3313.3157
$ws.Range("I62").Value = 2957.1428
$ws.Range("J62").Value = 3521.0833
$ws.Range("K62").Value = 2957.1428
$ws.Range("L62").Value = 3521.0833
$ws.Range("M62").Value = -2333.1428
$ws.Range("N62").Value = -4769.0833

$ws.Range("H65").Value = 3313.3157
$ws.Range("I65").Value = 2957.1428
$ws.Range("J65").Value = 3521.0833
$ws.Range("K65").Value = 14785.714
$ws.Range("L65").Value = 17605.4165
$ws.Range("M65").Value = -11665.714
$ws.Range("N65").Value = -23845.4165

$ws.Range("H132").Value = 1939.8823
$ws.Range("I132").Value = 728.6087
$ws.Range("J132").Value = 4472.5454
$ws.Range("K132").Value = 2185.8261
$ws.Range("L132").Value = 13417.6362
$ws.Range("M132").Value = 344.1738999999998
$ws.Range("N132").Value = -18477.6362

$ws.Range("H136").Value = 5849.7017
$ws.Range("I136").Value = 5046.8965
$ws.Range("J136").Value = 6681.1787
$ws.Range("K136").Value = 15140.6895
$ws.Range("L136").Value = 20043.5361
$ws.Range("M136").Value = -12590.6895
$ws.Range("N136").Value = -25143.5361
